$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.192.60"
$ws.Range("E2").Value = "  +5.38%  "
$ws.Range("D3").Value = "2.759.98"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'578.71"
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("D6").Value = "'154.30"
$ws.Range("E6").Value = "  +6.31%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").Value = "2.756.83"
$ws.Range("E9").Value = "  +3.06%  "
$ws.Range("E10").Value = "  +1.98%  "
$ws.Range("E11").Value = "  +4.85%  "
$ws.Range("E12").Value = "  +4.68%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "3.239.16"
$ws.Range("E14").Value = "  +2.98%  "
$ws.Range("D15").Value = "'26.29"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").Value = "64.053.60"
$ws.Range("E16").Value = "  +5.16%  "
$ws.Range("E17").Value = "  +6.12%  "
$ws.Range("D18").Value = "2.756.26"
$ws.Range("E18").Value = "  +2.97%  "
$ws.Range("D19").Value = "'11.94"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("E20").Value = "  +2.40%  "
$ws.Range("D21").Value = "'360.77"
$ws.Range("E21").Value = "  +2.65%  "
$ws.Range("D22").Value = "'6.97"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").Value = "'66.38"
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("E26").Value = "  +4.66%  "
$ws.Range("D27").Value = "'8.51"
$ws.Range("E27").Value = "  +4.51%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "0.0₃0917"
$ws.Range("E29").Value = "  +12.32%  "
$ws.Range("D30").Value = "'2.00"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").Value = "'7.10"
$ws.Range("E31").Value = "  +4.51%  "
$ws.Range("D32").Value = "'171.73"
$ws.Range("E32").Value = "  +3.09%  "
$ws.Range("D33").Value = "'1.23"
$ws.Range("E33").Value = "  +15.07%  "
$ws.Range("D35").Value = "'20.45"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").Value = "'4.79"
$ws.Range("E36").Value = "  +7.43%  "
$ws.Range("E37").Value = "  +8.16%  "
$ws.Range("B38").Value = "SuiNetwork"
$ws.Range("C38").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D38").Value = "'1.04"
$ws.Range("E38").Value = "  +18.25%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").Value = "'1.81"
$ws.Range("E39").Value = "  +9.73%  "
$ws.Range("D40").Value = "'345.93"
$ws.Range("E40").Value = "  +4.90%  "
$ws.Range("E41").Value = "  +5.15%  "
$ws.Range("E42").Value = "  +2.80%  "
$ws.Range("E43").Value = "  +7.63%  "
$ws.Range("E44").Value = "  +5.76%  "
$ws.Range("D45").Value = "'0.0592"
$ws.Range("E45").Value = "  +5.38%  "
$ws.Range("D46").Value = "'21.71"
$ws.Range("E46").Value = "  +5.82%  "
$ws.Range("E47").Value = "  +5.20%  "
$ws.Range("D48").Value = "'139.22"
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("E49").Value = "  +2.89%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("E51").Value = "  +0.01%  "
